$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = [double]$cell.Value2
        $cell.Value2 = [Math]::Floor($v + 0.5)
    }
}
